{"js": "// Change the report title from\n//   \"Coursera Capstone: Best location for City Waste Management in Kochi city\"\n// to\n//   \"Coursera Capstone: Determining best location for City Waste Management in Kochi city using Machine Learning\"\n// The \"Coursera Capstone: \" prefix keeps its original (bold, underlined,\n// Georgia 14pt) formatting; only the remainder of the title text changes,\n// and the replacement text inherits that same run formatting.\n\nconst body = context.document.body;\n\nconst oldTail = \"Best location for City Waste Management in Kochi city\";\nconst newTail = \"Determining best location for City Waste Management in Kochi city using Machine Learning\";\n\nconst results = body.search(oldTail, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find title text to update: \" + oldTail);\n}\n\n// Replace in place so the surrounding run formatting (bold/underline/font/\n// size/color) carries over to the new text automatically.\nresults.items[0].insertText(newTail, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Change the report title from\n#   \"Coursera Capstone: Best location for City Waste Management in Kochi city\"\n# to\n#   \"Coursera Capstone: Determining best location for City Waste Management in Kochi city using Machine Learning\"\n# The \"Coursera Capstone: \" prefix keeps its original run formatting\n# (bold, underlined, Georgia 14pt); only the remainder of the title text\n# is located and replaced, inheriting that same run formatting.\n\n$d = $word.ActiveDocument\n\n$oldTail = \"Best location for City Waste Management in Kochi city\"\n$newTail = \"Determining best location for City Waste Management in Kochi city using Machine Learning\"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $oldTail\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Forward = $true\n$find.Wrap = 1\n\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find title text to update: $oldTail\"\n}\n\n# $range now spans exactly the matched text (\"Execute\" without a\n# replacement collapses the range onto the hit); overwrite it in place so\n# the surrounding run formatting carries over to the new text.\n$range.Text = $newTail\n"}
